$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the imsi1 / imsi2 columns (I:J) - device stock & stock assignment merge
$ws.Range("I1:J1").EntireColumn.Delete()

# Update the active selection to match the new layout
$ws.Range("G7").Select()
